$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.306.21'
$ws.Range('E2').Value = '  +1.21%  '
$ws.Range('D3').Value = '1.853.63'
$ws.Range('E3').Value = '  +1.44%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  -0.89%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '314.11'
$ws.Range('E5').Value = '  +0.75%  '
$ws.Range('E6').Value = '  -0.74%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4617'
$ws.Range('E7').Value = '  -0.57%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3710'
$ws.Range('E8').Value = '  +0.09%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07309'
$ws.Range('E9').Value = '  -0.92%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8828'
$ws.Range('E10').Value = '  +1.08%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '19.90'
$ws.Range('E11').Value = '  -0.19%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07786'
$ws.Range('E12').Value = '  -0.40%  '
$ws.Range('D13').Value = '1.849.94'
$ws.Range('E13').Value = '  +3.79%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.375'
$ws.Range('E14').Value = '  +0.56%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.541'
$ws.Range('E15').Value = '  -0.77%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '91.76'
$ws.Range('E16').Value = '  +0.09%  '
$ws.Range('E17').Value = '  -0.91%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000009032'
$ws.Range('E18').Value = '  +2.16%  '
$ws.Range('E19').Value = '  -0.70%  '
$ws.Range('E20').Value = '  +0.88%  '
$ws.Range('D21').Value = '27.322.12'
$ws.Range('E21').Value = '  +1.31%  '
$ws.Range('E22').Value = '  -0.45%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.51'
$ws.Range('E23').Value = '  -0.72%  '
$ws.Range('D24').Value = '2.107.14'
$ws.Range('E24').Value = '  +4.15%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.925'
$ws.Range('E25').Value = '  +5.12%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '151.90'
$ws.Range('E26').Value = '  -0.37%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.41'
$ws.Range('E27').Value = '  +0.69%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.067'
$ws.Range('E28').Value = '  -0.63%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.097'
$ws.Range('E29').Value = '  -0.42%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '116.15'
$ws.Range('E30').Value = '  +0.45%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08827'
$ws.Range('E31').Value = '  -0.41%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7754'
$ws.Range('E32').Value = '  +6.65%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.084'
$ws.Range('E33').Value = '  +3.66%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.171'
$ws.Range('E34').Value = '  +3.12%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.499'
$ws.Range('E35').Value = '  +1.24%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.656'
$ws.Range('E36').Value = '  +6.53%  '
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.078'
$ws.Range('E37').Value = '  +0.34%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01958'
$ws.Range('E38').Value = '  +0.46%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.970'
$ws.Range('E39').Value = '  +1.45%  '
$ws.Range('B40').Value = 'Hedera'
$ws.Range('C40').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.05222'
$ws.Range('E40').Value = '  +0.02%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.987'
$ws.Range('E41').Value = '  -3.12%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5144'
$ws.Range('E42').Value = '  -0.91%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1632'
$ws.Range('E43').Value = '  +0.20%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.390'
$ws.Range('E44').Value = '  +2.12%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4811'
$ws.Range('E45').Value = '  -0.53%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.36'
$ws.Range('E46').Value = '  +1.77%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.001'
$ws.Range('E47').Value = '  -0.81%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '102.82'
$ws.Range('E48').Value = '  +0.25%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.646'
$ws.Range('E49').Value = '  +1.31%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06225'
$ws.Range('E50').Value = '  +0.13%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '65.71'
$ws.Range('E51').Value = '  +2.32%  '
